$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary  (row 2 values)
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.0498220640569395
$wsSummary.Range("C2").Value = 0.0498220640569395
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.09491525423728814
$wsSummary.Range("F2").Value = 0.2077151335311573
$wsSummary.Range("G2").Value = 0.5768621236133122
$wsSummary.Range("H2").Value = 0.7547819689673623
$wsSummary.Range("I2").Value = 28
$wsSummary.Range("J2").Value = 534
$wsSummary.Range("K2").Value = 0
$wsSummary.Range("L2").Value = 0

# ---------------------------------------------------------------------------
# Sheet: Classification Report
# ---------------------------------------------------------------------------
$wsReport = $wb.Worksheets.Item("Classification Report")

# Row 2 - class "0"
$wsReport.Range("B2").Value = 0
$wsReport.Range("C2").Value = 0
$wsReport.Range("D2").Value = 0

# Row 3 - class "1"
$wsReport.Range("B3").Value = 0.0498220640569395
$wsReport.Range("C3").Value = 1
$wsReport.Range("D3").Value = 0.09491525423728814

# Row 4 - accuracy
$wsReport.Range("B4").Value = 0.0498220640569395
$wsReport.Range("C4").Value = 0.0498220640569395
$wsReport.Range("D4").Value = 0.0498220640569395
$wsReport.Range("E4").Value = 0.0498220640569395

# Row 5 - macro avg
$wsReport.Range("B5").Value = 0.02491103202846975
$wsReport.Range("C5").Value = 0.5
$wsReport.Range("D5").Value = 0.04745762711864407

# Row 6 - weighted avg
$wsReport.Range("B6").Value = 0.002482238066893783
$wsReport.Range("C6").Value = 0.0498220640569395
$wsReport.Range("D6").Value = 0.004728873876590867

# ---------------------------------------------------------------------------
# Sheet: Confusion Matrix
# ---------------------------------------------------------------------------
$wsMatrix = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 - Actual 0
$wsMatrix.Range("B2").Value = 0
$wsMatrix.Range("C2").Value = 534

# Row 3 - Actual 1
$wsMatrix.Range("B3").Value = 0
$wsMatrix.Range("C3").Value = 28
